# Update LR-pairs data with new TPM-derived values (Cxcl1-Cxcr1).
# The table grows from 4 data rows (one per Sending cluster) to 8 data rows
# (one per Sending cluster x Target cluster [ECs, Resolving-Mac]), and all
# numeric columns are recomputed against the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl1"
$ws.Range("C2").Value = "Cxcr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.268639
$ws.Range("H2").Value = 60.80591700000001
$ws.Range("I2").Value = 0.09537690375401175
$ws.Range("J2").Value = 0.09537690375401174
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.0003833333333333333
$ws.Range("N2").Value = 0.00115
$ws.Range("O2").Value = 0.005317132262509131
$ws.Range("P2").Value = 0.005317132262509131
$ws.Range("Q2").Value = 0.007769644950000001
$ws.Range("R2").Value = 0.06992680455000001
$ws.Range("S2").Value = 0.0005071316120486842
$ws.Range("T2").Value = 0.0005071316120486841

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl1"
$ws.Range("C3").Value = "Cxcr1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.268639
$ws.Range("H3").Value = 60.80591700000001
$ws.Range("I3").Value = 0.09537690375401175
$ws.Range("J3").Value = 0.09537690375401174
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.07171066666666667
$ws.Range("N3").Value = 0.215132
$ws.Range("O3").Value = 0.9946828677374909
$ws.Range("P3").Value = 0.9946828677374908
$ws.Range("Q3").Value = 1.453477615116
$ws.Range("R3").Value = 13.081298536044
$ws.Range("S3").Value = 0.09486977214196307
$ws.Range("T3").Value = 0.09486977214196306

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl1"
$ws.Range("C4").Value = "Cxcr1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 127.7410203333333
$ws.Range("H4").Value = 383.223061
$ws.Range("I4").Value = 0.6011031624655011
$ws.Range("J4").Value = 0.601103162465501
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.0003833333333333333
$ws.Range("N4").Value = 0.00115
$ws.Range("O4").Value = 0.005317132262509131
$ws.Range("P4").Value = 0.005317132262509131
$ws.Range("Q4").Value = 0.04896739112777778
$ws.Range("R4").Value = 0.44070652015
$ws.Range("S4").Value = 0.003196145018241584
$ws.Range("T4").Value = 0.003196145018241583

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl1"
$ws.Range("C5").Value = "Cxcr1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 127.7410203333333
$ws.Range("H5").Value = 383.223061
$ws.Range("I5").Value = 0.6011031624655011
$ws.Range("J5").Value = 0.601103162465501
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.07171066666666667
$ws.Range("N5").Value = 0.215132
$ws.Range("O5").Value = 0.9946828677374909
$ws.Range("P5").Value = 0.9946828677374908
$ws.Range("Q5").Value = 9.160393728783557
$ws.Range("R5").Value = 82.44354355905202
$ws.Range("S5").Value = 0.5979070174472596
$ws.Range("T5").Value = 0.5979070174472594

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cxcl1"
$ws.Range("C6").Value = "Cxcr1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 42.02733833333333
$ws.Range("H6").Value = 126.082015
$ws.Range("I6").Value = 0.1977654939365007
$ws.Range("J6").Value = 0.1977654939365007
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.0003833333333333333
$ws.Range("N6").Value = 0.00115
$ws.Range("O6").Value = 0.005317132262509131
$ws.Range("P6").Value = 0.005317132262509131
$ws.Range("Q6").Value = 0.01611047969444444
$ws.Range("R6").Value = 0.14499431725
$ws.Range("S6").Value = 0.001051545288220822
$ws.Range("T6").Value = 0.001051545288220822

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cxcl1"
$ws.Range("C7").Value = "Cxcr1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 42.02733833333333
$ws.Range("H7").Value = 126.082015
$ws.Range("I7").Value = 0.1977654939365007
$ws.Range("J7").Value = 0.1977654939365007
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.07171066666666667
$ws.Range("N7").Value = 0.215132
$ws.Range("O7").Value = 0.9946828677374909
$ws.Range("P7").Value = 0.9946828677374908
$ws.Range("Q7").Value = 3.013808450108888
$ws.Range("R7").Value = 27.12427605098
$ws.Range("S7").Value = 0.1967139486482798
$ws.Range("T7").Value = 0.1967139486482798

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Cxcl1"
$ws.Range("C8").Value = "Cxcr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 22.47397933333333
$ws.Range("H8").Value = 67.421938
$ws.Range("I8").Value = 0.1057544398439867
$ws.Range("J8").Value = 0.1057544398439867
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.0003833333333333333
$ws.Range("N8").Value = 0.00115
$ws.Range("O8").Value = 0.005317132262509131
$ws.Range("P8").Value = 0.005317132262509131
$ws.Range("Q8").Value = 0.008615025411111111
$ws.Range("R8").Value = 0.0775352287
$ws.Range("S8").Value = 0.0005623103439980427
$ws.Range("T8").Value = 0.0005623103439980426

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Cxcl1"
$ws.Range("C9").Value = "Cxcr1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 22.47397933333333
$ws.Range("H9").Value = 67.421938
$ws.Range("I9").Value = 0.1057544398439867
$ws.Range("J9").Value = 0.1057544398439867
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.07171066666666667
$ws.Range("N9").Value = 0.215132
$ws.Range("O9").Value = 0.9946828677374909
$ws.Range("P9").Value = 0.9946828677374908
$ws.Range("Q9").Value = 1.611624040646222
$ws.Range("R9").Value = 14.504616365816
$ws.Range("S9").Value = 0.1051921294999886
$ws.Range("T9").Value = 0.1051921294999886
